$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: Sheet0 -> Sheet1
$ws.Name = "Sheet1"

# Update rows 8-11 (A:D) with new values per the commit's data refresh
$ws.Range("A8").Value = "pepsi"
$ws.Range("B8").Value = 2.1
$ws.Range("C8").Value = "JE"
$ws.Range("D8").Value = "drink"

$ws.Range("A9").Value = "baks0"
$ws.Range("B9").Value = 7.0
$ws.Range("C9").Value = "NTU"
$ws.Range("D9").Value = "BURGER"

$ws.Range("A10").Value = "pizza"
$ws.Range("B10").Value = 50.0
$ws.Range("C10").Value = "NTU"
$ws.Range("D10").Value = "SET MEAL"

$ws.Range("A11").Value = "Sprite"
$ws.Range("B11").Value = 6.0
$ws.Range("C11").Value = "NTU"
$ws.Range("D11").Value = "BURGER"
